$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$group = "Srv_Geo_Bk"
$user = "ansible_test"

$newRows = @(
    @("172.16.11.7",   "akb68-n1"),
    @("172.16.3.6",    "akvil-n1"),
    @("172.16.107.7",  "akbar-bk"),
    @("172.16.116.80", "aksin-bk"),
    @("10.16.77.4",    "akcan-n1"),
    @("10.173.3.7",    "alapa-bk"),
    @("10.121.0.6",    "akyop-n1"),
    @("172.16.134.4",  "akede-n1"),
    @("10.125.3.8",    "albu2-bk"),
    @("10.245.3.4",    "ktnqs-n1")
)

$startRow = 119
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ip = $newRows[$i][0]
    $host_ = $newRows[$i][1]

    $ws.Cells.Item($row, 1).Value = $ip
    $ws.Cells.Item($row, 2).Value = $group
    $ws.Cells.Item($row, 3).Value = $ip
    $ws.Cells.Item($row, 4).Value = $user
    $ws.Cells.Item($row, 5).Value = $host_
}

$ws.Range("A119").Select()
